$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2:E2) - CON
$ws.Range("B2").Value = 64.409404315482547
$ws.Range("C2").Value = 34.377179995401043
$ws.Range("D2").Value = 63.043111677046504
$ws.Range("E2").Value = 38.247644724786852

# Row 3 data values (B3:E3) - STR
$ws.Range("B3").Value = 61.040297496419768
$ws.Range("C3").Value = 42.45096915661842
$ws.Range("D3").Value = 50.170624494490227
$ws.Range("E3").Value = 44.281495040670407

# Update the active selection to match the updated data range
$ws.Range("B1:E3").Select()
